$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 1283; this shifts rows 1283:1528 down to 1284:1529
$ws.Rows.Item(1283).Insert()

# Populate the newly inserted row 1283 with the new record's data
$ws.Cells.Item(1283, 1).Value = 756
$ws.Cells.Item(1283, 2).Value = "RECURSOS DE ALIENAÇÃO DE BENS/ATIVOS - ADMINISTRAÇÃO INDIRETA"
$ws.Cells.Item(1283, 3).Value = 2251
$ws.Cells.Item(1283, 4).Value = "JUCEMG"
$ws.Cells.Item(1283, 5).Value = 47
$ws.Cells.Item(1283, 6).Value = "ALIENACAO DE BENS DE ENTIDADES ESTADUAIS"
$ws.Cells.Item(1283, 7).Value = 2213010101000
$ws.Cells.Item(1283, 8).Value = "ALIENACAO BENS MOVEIS SEMOVENTES - PRINC."
